$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String literals (kept as variables for clarity/reuse) ---
$objetivos_pt = 'Introduzir conceitos básicos sobre a influência da estrutura química dos materiais poliméricos na morfologia e nas propriedades mecânicas, ópticas, elétricas e químicas.'
$docentes = '5840897 - Clodoaldo Saron'
$programa_resumido_pt = 'Histórico dos polímeros, Estrutura e Propriedades de Materiais Poliméricos, Polímeros Termoplásticos, Elastômeros, Polímeros Termorrígidos e Aditivos para Polímeros'
$programa_pt = 'Introdução: nomenclatura, arquitetura molecular e estrutura configuracional. Estado sólido: amorfo, cristalino e elastomérico. Estrutura e propriedades. Thermoplásticos: estrutura, propriedades e aplicações.  Elastômeros: estrutura, propriedades e aplicações. Resinas termorrígidas: estrutura, propriedades e aplicações.  Aditivos para polímeros: classes e aplicações. Avaliação'
$metodo = 'Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula'
$criterio = 'Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação.'
$norma_recuperacao = 'O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5.'
$bibliografia = '1. L. C. SAWYER & D. T. GRUBB. Polymer Microscopy. New York: Chapman and Hall, 1987.2. H. G. ELLIAS Macromolecules -1-Structure and Properties, 2nd ed. Vol.1, New York, Plenum Press, 1984.3. J. M. G. COWIE Polymers: Chemistry and Physics of Modern Materials", New York: Chapman & Hall, 1998.4. C. A. HARPER Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992.5. S. V. CANEVAROLO Jr. Ciência dos Polímeros. Editora Artiliber, 2001.6. E. B. MANO & L.C. MENDES Introdução a Polímeros. São Paulo: Ed. Edgard Blücher, 1999.7. E. B. MANO Polímeros como Materiais de Engenharia. São Paulo: Ed. Edgard Blücher, 1991.8. E. B. MANO & L. C. MENDES Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000.'

# --- Step 1: insert a new row at 13 (shifts old rows 13-23 down to 14-24) ---
$ws.Rows("13").Insert()

# The inserted row mistakenly carries a placeholder A13 cell (style only, no value);
# the target layout has no A13 cell on the new row, so remove it entirely.
$ws.Range("A13").Clear()

# --- Step 2: populate the new row 13 (Docentes responsaveis value), using
# Copy/PasteSpecial(formats) from row 10 (still B=style2/C=style3) so the exact
# existing cellXfs indices are reused instead of synthesizing new style combos. ---
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B13").Value = $docentes
$ws.Range("C13").Value = $docentes

# --- Step 3: content updates (rows identified by their FINAL row number, i.e. after the insert/shift) ---

# Row 10 "Objetivos:" -- replace the misplaced "5840897 - Clodoaldo Saron" with the real Portuguese objective text
$ws.Range("B10").Value = $objetivos_pt
$ws.Range("C10").Value = $objetivos_pt

# Row 14 "Programa resumido:" -- replace "Semestral" with the real summary syllabus text
$ws.Range("B14").Value = $programa_resumido_pt
$ws.Range("C14").Value = $programa_resumido_pt

# Row 16 "Programa:" -- replace the stray date with the real full syllabus text
$ws.Range("B16").Value = $programa_pt
$ws.Range("C16").Value = $programa_pt

# Row 19 "Metodo:" -- replace the misplaced docentes value with the real method text
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Row 20 "Criterio:" -- replace the misplaced method text with the real criterio text
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Row 21 "Norma de recuperacao:" -- replace the misplaced criterio text with the real recovery-norm text
$ws.Range("B21").Value = $norma_recuperacao
$ws.Range("C21").Value = $norma_recuperacao

# Row 22 "Bibliografia:" -- replace the misplaced recovery-norm text with the real bibliography text
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

Write-Host "done"
